$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:H14").Value = -0.0005817

$ws.Columns("A").AutoFit() | Out-Null

$ws.Range("I7").Select()
